# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# The W/L/T values are constant for every player row since they
# represent the team's season record, not a per-player stat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the header text first so the paste-special below only needs
# to bring over formatting (bold font, thin border, centered/top
# alignment) that matches the rest of row 1's header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the last existing header cell (AC1) onto
# the three new header cells so they share the same style as the
# other headers in row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Re-apply the header labels (PasteSpecial only copied formats, but
# set again defensively in case of any clobbering).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every player row (2-52) with the
# season's win/loss/tie totals.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
